$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows([int]$r1, [int]$r2) {
    # Columns B, and E through AD (inclusive) are swapped between the two rows.
    # Column A (row index), C, D remain untouched.
    $cols = @('B') + @('E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z','AA','AB','AC','AD')
    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

Swap-Rows 42 43
Swap-Rows 139 140
Swap-Rows 212 213
Swap-Rows 301 302
